# "Finish add customer step"
# Inserts two new worksheets - "customer" and "assignProjectTest" - between
# the existing "addLeaveTypeTest" and "employee" sheets, and populates them
# with their testcase data. The pre-existing "employee" and "addSkillTest"
# sheets are left untouched (they simply shift two positions to the right).

$wb = $excel.ActiveWorkbook

# Anchor the new sheets right after "addLeaveTypeTest" (before "employee").
$afterLeaveType = $wb.Worksheets.Item("addLeaveTypeTest")

$customer = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterLeaveType)
$customer.Name = "customer"

$assignProject = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $customer)
$assignProject.Name = "assignProjectTest"

# ---- customer sheet -------------------------------------------------
$customer.Range("A1").Value = "Testcase ID"
$customer.Range("B1").Value = "Index in testcase"
$customer.Range("C1").Value = "Name"
$customer.Range("D1").Value = "Description"

$customer.Range("A2").Value = "ASSIGN_PROJECT_TC000"
$customer.Range("B2").Value = 0
$customer.Range("C2").Value = "John"
$customer.Range("D2").Value = "test"

$customer.Range("A3").Value = "ASSIGN_PROJECT_TC000"
$customer.Range("B3").Value = 1
$customer.Range("C3").Value = "Jame"
$customer.Range("D3").Value = "test"

$customer.Columns.Item(1).ColumnWidth = 27.92
$customer.Columns.Item(2).ColumnWidth = 14.58

# ---- assignProjectTest sheet -----------------------------------------
$assignProject.Range("A1").Value = "Testcase Id"
$assignProject.Range("A2").Value = "ASSIGN_PROJECT_TC000"

$assignProject.Columns.Item(1).ColumnWidth = 26.75

# "customer" is the tab that ends up selected/active, matching the
# original workbook where the sheet inserted at this slot (previously
# "employee") held the active tab.
$customer.Activate()
$customer.Range("A1").Select()
